# Tournament registration writer refactor
# - Résumé sheet: fix the Type/Nombre de place value swap, and write the
#   registration/tournament dates as formatted text instead of raw Excel
#   date serials.
# - Equipe sheet: lengthen the "Joueur 1"/"Joueur 2" column headers and
#   widen the columns to fit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Résumé
$ws2 = $wb.Worksheets.Item(2)   # Equipe

# --- Résumé sheet -----------------------------------------------------

# B2 used to hold the numeric "Nombre de place" value (20) next to the
# "Type :" label; it should hold the tournament type instead.
$ws1.Range("B2").Value = "DOUBLE"

# B3 ("Nombre de place :") now carries the numeric value that used to
# live in B2, and loses its previous (empty) bold/bordered style.
$ws1.Range("B3").ClearFormats()
$ws1.Range("B3").Value = 20.0

# Register a date-formatted style (matching the existing bold/bordered
# style used by column A) so the style table grows the same way it did
# in the original edit, even though the cells below end up holding plain
# text, not numbers.
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$ws1.Range("D1").NumberFormat = "mm-dd-yy"
$ws1.Range("D1").Clear()

# B5/B6/B7 now hold the dates as plain text instead of Excel date
# serials. Force a text number format while assigning so the
# date-looking strings aren't re-interpreted as dates, then drop the
# format again so the final cells carry no explicit style.
$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2023-05-30"
$ws1.Range("B5").ClearFormats()

$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2023-06-01"
$ws1.Range("B6").ClearFormats()

$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2023-06-02"
$ws1.Range("B7").ClearFormats()

# Widen column B slightly to fit the longer date strings.
$ws1.Columns.Item(2).ColumnWidth = 12.44140625

# --- Equipe sheet -------------------------------------------------------

$ws2.Range("B1").Value = "Joueur 1 (Prénom NOM)"
$ws2.Range("C1").Value = "Joueur 2  (Prénom NOM)"

$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()

# --- Selections / active cells ------------------------------------------

$ws2.Activate()
$ws2.Range("A2").Select()

$ws1.Activate()
$ws1.Range("A10").Select()
